$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Actual" header to "ActualValue"
$ws.Range("C2").Value = "ActualValue"

# Add a new "ActualPage" column in D, copying C2's formatting (fill/border/font)
# so the new cell gets the same style index as C2 without creating a new style.
$ws.Range("C2").Copy($ws.Range("D2"))
$ws.Range("D2").Value = "ActualPage"

# Give column D a sensible custom width (matches the ~13.3 char width used for
# the other header columns in this template).
$ws.Columns.Item(4).ColumnWidth = 12.5
